# Apply cryptos list update (prices & volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.925.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.64"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5128"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2545"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06168"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("E10").Value = "  -4.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07507"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.635.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.331"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.849.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5344"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7905"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.61"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.904.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.600"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "183.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.884"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.021"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1187"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.257"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.370"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05919"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.239"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.335"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.317"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.591"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9600"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("E36").Value = "  -1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5753"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01574"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.051.08"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.784.68"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05194"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.861"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4223"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.68%  "

# Rows 42 and 43: Coin/Link swapped (TrustWalletToken <-> FraxShare) with updated Price/Volume
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.740"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.34%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8350"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.75%  "
